# Refresh the "想去人数" (interest count) figures in column F across the
# 展览, 演出 and 全部类型 sheets to match the newly scraped snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1677
$ws1.Range("F3").Value = 9114
$ws1.Range("F4").Value = 114
$ws1.Range("F8").Value = 203
$ws1.Range("F10").Value = 96
$ws1.Range("F11").Value = 5911
$ws1.Range("F15").Value = 4452
$ws1.Range("F18").Value = 1150
$ws1.Range("F20").Value = 339
$ws1.Range("F22").Value = 259
$ws1.Range("F24").Value = 2758

# 演出 (sheet2): row -> new F value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 35

# 全部类型 (sheet4): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1677
$ws4.Range("F3").Value = 9114
$ws4.Range("F4").Value = 114
$ws4.Range("F5").Value = 35
$ws4.Range("F9").Value = 203
$ws4.Range("F11").Value = 96
$ws4.Range("F12").Value = 5911
$ws4.Range("F16").Value = 4452
$ws4.Range("F19").Value = 1150
$ws4.Range("F21").Value = 339
$ws4.Range("F23").Value = 259
$ws4.Range("F25").Value = 2758
